$d = $word.ActiveDocument

# Replace the first occurrence: "Venus orbiter radiometric temperature experiment (VORTEX"
# becomes "Venus orbiter infrared radiometer (OIR"
$d.Content.Find.Execute(
    "Venus orbiter radiometric temperature experiment (VORTEX",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Venus orbiter infrared radiometer (OIR",
    2)

# Replace the second occurrence: ", Taylor et al., 1979).  VORTEX used the same principles"
# becomes ", Taylor et al., 1979).  OIR used the same principles"
$d.Content.Find.Execute(
    ", Taylor et al., 1979).  VORTEX used the same principles",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    ", Taylor et al., 1979).  OIR used the same principles",
    2)
